{"js": "// Remove the standalone \"Meta description\" paragraph (right after the H1 title)\n// and relocate its gist to the end of the document, splitting the final\n// paragraph (the image-prompt paragraph) into a new bold title paragraph\n// followed by the (now re-worded) italic paragraph.\n\nconst body = context.document.body;\n\n// --- Step 1: delete the \"Meta description: ...\" paragraph -----------------\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet metaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Meta description\") === 0) {\n    metaIndex = i;\n    break;\n  }\n}\nif (metaIndex === -1) {\n  throw new Error('Could not find the \"Meta description\" paragraph.');\n}\nparagraphs.items[metaIndex].delete();\nawait context.sync();\n\n// --- Step 2: split the last paragraph into two -----------------------------\nparagraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst wholeRange = lastParagraph.getRange(Word.RangeLocation.whole);\n\n// Replace the whole (currently-italic, \"For the feature image...\" ) paragraph\n// with two paragraphs: a new bold title paragraph, and the original italic\n// paragraph re-worded to hold the ex-meta-description text.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +\n  '<w:t>Play Delicious Candy PopWins for Free: Review &amp; Features</w:t>' +\n  '</w:r></w:p>' +\n  '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr>' +\n  '<w:t>Read our review of Delicious Candy PopWins, a slot game featuring expanding grids and free spins. Play for free and discover the features and mechanics.</w:t>' +\n  '</w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nwholeRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Remove the standalone \"Meta description\" paragraph (right after the H1\n# title) and relocate its gist to the end of the document, splitting the\n# final paragraph (the image-prompt paragraph) into a new bold title\n# paragraph followed by the (now re-worded) italic paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"Meta description: ...\" paragraph ------------------\n$metaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"Meta description\")) {\n        $metaIndex = $i\n        break\n    }\n}\nif ($metaIndex -eq -1) {\n    throw \"Could not find the 'Meta description' paragraph.\"\n}\n$d.Paragraphs.Item($metaIndex).Range.Delete()\n\n# --- Step 2: split the last paragraph into two ------------------------------\n$count = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($count)\n$rng = $lastParagraph.Range\n\n# Exclude the trailing paragraph mark so the inserted content doesn't leave\n# a stray empty paragraph behind.\n[void]$rng.MoveEnd(1, -1)\n\n# Replace the whole (currently-italic, \"For the feature image...\" ) paragraph\n# text with two paragraphs: a new bold title paragraph, and the original\n# italic paragraph re-worded to hold the ex-meta-description text.\n$xml = \"<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Delicious Candy PopWins for Free: Review &amp; Features</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Delicious Candy PopWins, a slot game featuring expanding grids and free spins. Play for free and discover the features and mechanics.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n\n$rng.InsertXML($xml)\n"}
